$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 14.3280885520863
$ws.Range("D2").Value2 = 3.357149727881108
$ws.Range("E2").Value2 = 20.73015233425628
$ws.Range("F2").Value2 = 21.64544918826361
$ws.Range("G2").Value2 = 3.587041825821306
$ws.Range("I2").Value2 = 23.80554090927901
$ws.Range("L2").Value2 = 9.126775261208101
$ws.Range("M2").Value2 = 14.46051191956086
$ws.Range("N2").Value2 = 18.67389829705905
$ws.Range("O2").Value2 = 18.84889308477368
$ws.Range("B3").Value2 = 14.08723468803584
$ws.Range("D3").Value2 = 3.345270998204921
$ws.Range("E3").Value2 = 20.54091440889676
$ws.Range("F3").Value2 = 21.36950019760859
$ws.Range("G3").Value2 = 3.589551984262742
$ws.Range("I3").Value2 = 23.9501469762957
$ws.Range("L3").Value2 = 9.043663792821205
$ws.Range("M3").Value2 = 14.35844328079726
$ws.Range("N3").Value2 = 18.68594310877313
$ws.Range("O3").Value2 = 18.71455813311098
$ws.Range("B4").Value2 = 13.93899633274005
$ws.Range("D4").Value2 = 3.337859394778944
$ws.Range("E4").Value2 = 20.42377429749029
$ws.Range("F4").Value2 = 21.20569624594085
$ws.Range("G4").Value2 = 3.591176289876296
$ws.Range("I4").Value2 = 24.04344051476142
$ws.Range("L4").Value2 = 8.993569860024715
$ws.Range("M4").Value2 = 14.29718447948174
$ws.Range("N4").Value2 = 18.69545339925583
$ws.Range("O4").Value2 = 18.63752583019081
$ws.Range("B5").Value2 = 13.8785720011147
$ws.Range("D5").Value2 = 3.334810120241796
$ws.Range("E5").Value2 = 20.37583068868311
$ws.Range("F5").Value2 = 21.14044714153372
$ws.Range("G5").Value2 = 3.591859162399067
$ws.Range("I5").Value2 = 24.08259445977117
$ws.Range("L5").Value2 = 8.973410298413159
$ws.Range("M5").Value2 = 14.27259583260069
$ws.Range("N5").Value2 = 18.69986238344344
$ws.Range("O5").Value2 = 18.60753515820188
$ws.Range("B6").Value2 = 13.86853976868662
$ws.Range("D6").Value2 = 3.334302059957579
$ws.Range("E6").Value2 = 20.36785800408807
$ws.Range("F6").Value2 = 21.12970577779968
$ws.Range("G6").Value2 = 3.591973820513253
$ws.Range("I6").Value2 = 24.08916465533205
$ws.Range("L6").Value2 = 8.970078700013351
$ws.Range("M6").Value2 = 14.26853611186678
$ws.Range("N6").Value2 = 18.70062676008625
$ws.Range("O6").Value2 = 18.60264065143169
$ws.Range("B7").Value2 = 13.93818139690263
$ws.Range("D7").Value2 = 3.337818387659469
$ws.Range("E7").Value2 = 20.42312851576629
$ws.Range("F7").Value2 = 21.20481007895513
$ws.Range("G7").Value2 = 3.591185414397016
$ws.Range("I7").Value2 = 24.04396395427431
$ws.Range("L7").Value2 = 8.993296928040742
$ws.Range("M7").Value2 = 14.29685132506205
$ws.Range("N7").Value2 = 18.6955106981932
$ws.Range("O7").Value2 = 18.63711565667602
$ws.Range("B8").Value2 = 14.24515873157063
$ws.Range("D8").Value2 = 3.353078848783796
$ws.Range("E8").Value2 = 20.66511774427303
$ws.Range("F8").Value2 = 21.54918653082406
$ws.Range("G8").Value2 = 3.587890131397932
$ws.Range("I8").Value2 = 23.85446809365834
$ws.Range("L8").Value2 = 9.097935239844718
$ws.Range("M8").Value2 = 14.42503769707015
$ws.Range("N8").Value2 = 18.67761347322687
$ws.Range("O8").Value2 = 18.80146222137145
$ws.Range("B9").Value2 = 14.84120977821415
$ws.Range("D9").Value2 = 3.382038512161989
$ws.Range("E9").Value2 = 21.13082912781159
$ws.Range("F9").Value2 = 22.26521677557766
$ws.Range("G9").Value2 = 3.582083972934912
$ws.Range("I9").Value2 = 23.51845943645343
$ws.Range("L9").Value2 = 9.309731298711887
$ws.Range("M9").Value2 = 14.68671859091384
$ws.Range("N9").Value2 = 18.65922057221555
$ws.Range("O9").Value2 = 19.16553099925294
$ws.Range("B10").Value2 = 15.27145324369149
$ws.Range("D10").Value2 = 3.402687451931108
$ws.Range("E10").Value2 = 21.46586856217897
$ws.Range("F10").Value2 = 22.81064595417782
$ws.Range("G10").Value2 = 3.578213634806172
$ws.Range("I10").Value2 = 23.29307768986745
$ws.Range("L10").Value2 = 9.468227722184572
$ws.Range("M10").Value2 = 14.88403431863271
$ws.Range("N10").Value2 = 18.65578624993435
$ws.Range("O10").Value2 = 19.45636592988428
$ws.Range("B11").Value2 = 15.4647022735226
$ws.Range("D11").Value2 = 3.411936653978399
$ws.Range("E11").Value2 = 21.61637908918142
$ws.Range("F11").Value2 = 23.06180192738228
$ws.Range("G11").Value2 = 3.576537850137083
$ws.Range("I11").Value2 = 23.19516508569138
$ws.Range("L11").Value2 = 9.540715278338011
$ws.Range("M11").Value2 = 14.97462960923945
$ws.Range("N11").Value2 = 18.65638884116142
$ws.Range("O11").Value2 = 19.59325127790343
$ws.Range("B12").Value2 = 15.53746126403533
$ws.Range("D12").Value2 = 3.415417641901388
$ws.Range("E12").Value2 = 21.67307031968667
$ws.Range("F12").Value2 = 23.15724694089287
$ws.Range("G12").Value2 = 3.575915404393716
$ws.Range("I12").Value2 = 23.15874830844864
$ws.Range("L12").Value2 = 9.568199544941944
$ws.Range("M12").Value2 = 15.00903302504209
$ws.Range("N12").Value2 = 18.656926182407
$ws.Range("O12").Value2 = 19.64570140900254
$ws.Range("B13").Value2 = 15.52181100287785
$ws.Range("D13").Value2 = 3.414668916505787
$ws.Range("E13").Value2 = 21.66087484184045
$ws.Range("F13").Value2 = 23.13667769299729
$ws.Range("G13").Value2 = 3.576048920415339
$ws.Range("I13").Value2 = 23.16656197985854
$ws.Range("L13").Value2 = 9.562279098170535
$ws.Range("M13").Value2 = 15.00161969173253
$ws.Range("N13").Value2 = 18.65679674189227
$ws.Range("O13").Value2 = 19.63437871915953
$ws.Range("B14").Value2 = 15.47069698188423
$ws.Range("D14").Value2 = 3.412223471887067
$ws.Range("E14").Value2 = 21.62104937294264
$ws.Range("F14").Value2 = 23.06964804075385
$ws.Range("G14").Value2 = 3.576486398261874
$ws.Range("I14").Value2 = 23.19215583308842
$ws.Range("L14").Value2 = 9.542975846642204
$ws.Range("M14").Value2 = 14.97745818896085
$ws.Range("N14").Value2 = 18.65642686530122
$ws.Range("O14").Value2 = 19.59755432249101
$ws.Range("B15").Value2 = 15.43933157023868
$ws.Range("D15").Value2 = 3.41072274596989
$ws.Range("E15").Value2 = 21.59661465863373
$ws.Range("F15").Value2 = 23.02863150388424
$ws.Range("G15").Value2 = 3.576755944849083
$ws.Range("I15").Value2 = 23.20791874879971
$ws.Range("L15").Value2 = 9.531155964993662
$ws.Range("M15").Value2 = 14.96267050724826
$ws.Range("N15").Value2 = 18.65624050068331
$ws.Range("O15").Value2 = 19.5750770191942
$ws.Range("B16").Value2 = 15.25876862393064
$ws.Range("D16").Value2 = 3.402080029912465
$ws.Range("E16").Value2 = 21.45599147637835
$ws.Range("F16").Value2 = 22.79428494332398
$ws.Range("G16").Value2 = 3.578324853222418
$ws.Range("I16").Value2 = 23.29956911229858
$ws.Range("L16").Value2 = 9.463496489969446
$ws.Range("M16").Value2 = 14.87812853307056
$ws.Range("N16").Value2 = 18.655790249902
$ws.Range("O16").Value2 = 19.44750860492377
$ws.Range("B17").Value2 = 15.14731907623625
$ws.Range("D17").Value2 = 3.39674055254592
$ws.Range("E17").Value2 = 21.36921500645752
$ws.Range("F17").Value2 = 22.65122801514811
$ws.Range("G17").Value2 = 3.579309014518291
$ws.Range("I17").Value2 = 23.35697344082174
$ws.Range("L17").Value2 = 9.422073799216129
$ws.Range("M17").Value2 = 14.82646225985493
$ws.Range("N17").Value2 = 18.6560670077166
$ws.Range("O17").Value2 = 19.37039215247129
$ws.Range("B18").Value2 = 15.0829873627227
$ws.Range("D18").Value2 = 3.393655908152409
$ws.Range("E18").Value2 = 21.31912637467251
$ws.Range("F18").Value2 = 22.56923723151346
$ws.Range("G18").Value2 = 3.579883068491405
$ws.Range("I18").Value2 = 23.39042542486482
$ws.Range("L18").Value2 = 9.398286204162181
$ws.Range("M18").Value2 = 14.79682520947107
$ws.Range("N18").Value2 = 18.65643009457341
$ws.Range("O18").Value2 = 19.32647093796423
$ws.Range("B19").Value2 = 15.06116846853711
$ws.Range("D19").Value2 = 3.392609202033812
$ws.Range("E19").Value2 = 21.30213774743626
$ws.Range("F19").Value2 = 22.54152971356838
$ws.Range("G19").Value2 = 3.580078807625799
$ws.Range("I19").Value2 = 23.40182641464816
$ws.Range("L19").Value2 = 9.390239234381776
$ws.Range("M19").Value2 = 14.78680506278759
$ws.Range("N19").Value2 = 18.65658811436092
$ws.Range("O19").Value2 = 19.31167586436349
$ws.Range("B20").Value2 = 15.15920723941296
$ws.Range("D20").Value2 = 3.397310353628585
$ws.Range("E20").Value2 = 21.37847106231418
$ws.Range("F20").Value2 = 22.66642719336505
$ws.Range("G20").Value2 = 3.579203422354384
$ws.Range("I20").Value2 = 23.35081770125664
$ws.Range("L20").Value2 = 9.426479570838298
$ws.Range("M20").Value2 = 14.83195410843256
$ws.Range("N20").Value2 = 18.65601645850646
$ws.Range("O20").Value2 = 19.37855671755736
$ws.Range("B21").Value2 = 15.48572231518275
$ws.Range("D21").Value2 = 3.41294234739271
$ws.Range("E21").Value2 = 21.6327555649599
$ws.Range("F21").Value2 = 23.08932789221892
$ws.Range("G21").Value2 = 3.576357571512093
$ws.Range("I21").Value2 = 23.1846203931302
$ws.Range("L21").Value2 = 9.548644900074802
$ws.Range("M21").Value2 = 14.98455256188626
$ws.Range("N21").Value2 = 18.65652713428438
$ws.Range("O21").Value2 = 19.60835421968157
$ws.Range("B22").Value2 = 15.6966415739574
$ws.Range("D22").Value2 = 3.423033180109601
$ws.Range("E22").Value2 = 21.79716046469607
$ws.Range("F22").Value2 = 23.36764146679346
$ws.Range("G22").Value2 = 3.574568360982708
$ws.Range("I22").Value2 = 23.07984996589733
$ws.Range("L22").Value2 = 9.628680336701013
$ws.Range("M22").Value2 = 15.08483878121514
$ws.Range("N22").Value2 = 18.65866181186783
$ws.Range("O22").Value2 = 19.76210127684596
$ws.Range("B23").Value2 = 15.58431727920236
$ws.Range("D23").Value2 = 3.417659254784462
$ws.Range("E23").Value2 = 21.70958770280811
$ws.Range("F23").Value2 = 23.21895661654745
$ws.Range("G23").Value2 = 3.575516846686022
$ws.Range("I23").Value2 = 23.13541669938742
$ws.Range("L23").Value2 = 9.585953058002907
$ws.Range("M23").Value2 = 15.03127108415476
$ws.Range("N23").Value2 = 18.65735844335829
$ws.Range("O23").Value2 = 19.67973270377378
$ws.Range("B24").Value2 = 15.15383339909426
$ws.Range("D24").Value2 = 3.397052792921493
$ws.Range("E24").Value2 = 21.37428701920085
$ws.Range("F24").Value2 = 22.65955484435898
$ws.Range("G24").Value2 = 3.579251134878789
$ws.Range("I24").Value2 = 23.35359931080503
$ws.Range("L24").Value2 = 9.424487636209541
$ws.Range("M24").Value2 = 14.82947103451551
$ws.Range("N24").Value2 = 18.656038676274
$ws.Range("O24").Value2 = 19.37486422448575
$ws.Range("B25").Value2 = 14.68102491154187
$ws.Range("D25").Value2 = 3.374311575651968
$ws.Range("E25").Value2 = 21.005967478014
$ws.Range("F25").Value2 = 22.06771169624221
$ws.Range("G25").Value2 = 3.583584932229588
$ws.Range("I25").Value2 = 23.60557052457187
$ws.Range("L25").Value2 = 9.251844544381919
$ws.Range("M25").Value2 = 14.6149484974525
$ws.Range("N25").Value2 = 18.66241853287728
$ws.Range("O25").Value2 = 19.06277758319763
